$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 267.7857
$ws.Range("I4").Value = 211.38461
$ws.Range("J4").Value = 1001
$ws.Range("K4").Value = 211.38461
$ws.Range("L4").Value = 1001
$ws.Range("M4").Value = -97.38461000000001
$ws.Range("N4").Value = -1229

$ws.Range("H5").Value = 236.4
$ws.Range("I5").Value = 45.5
$ws.Range("K5").Value = 45.5
$ws.Range("M5").Value = 69.5

$ws.Range("H17").Value = 1676103.2
$ws.Range("J17").Value = 5026250
$ws.Range("L17").Value = 15078750
$ws.Range("N17").Value = -15079086

$ws.Range("H64").Value = 6727.5
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496

$ws.Range("H67").Value = 6727.5
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716

$ws.Range("H106").Value = 3342515.8
$ws.Range("I106").Value = 4175648
$ws.Range("K106").Value = 4175648
$ws.Range("M106").Value = -4175017

$ws.Range("H116").Value = 23923.416
$ws.Range("I116").Value = 25507.363
$ws.Range("J116").Value = 6500
$ws.Range("K116").Value = 25507.363
$ws.Range("L116").Value = 6500
$ws.Range("M116").Value = -22065.363
$ws.Range("N116").Value = -13384

$ws.Range("H132").Value = 7815.263
$ws.Range("I132").Value = 8868.125
$ws.Range("K132").Value = 26604.375
$ws.Range("M132").Value = -24074.375

$ws.Range("H136").Value = 125000
$ws.Range("J136").Value = 125000
$ws.Range("L136").Value = 125000
$ws.Range("N136").Value = -135200

$ws.Range("H137").Value = 412500
$ws.Range("I137").Value = 412500
$ws.Range("K137").Value = 1237500
$ws.Range("M137").Value = -1234950

$ws.Range("H141").Value = 4496.8335
$ws.Range("I141").Value = 4496.2
$ws.Range("K141").Value = 13488.6
$ws.Range("M141").Value = -8308.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 20002000
$ws.Range("J10").Value = 3999
$ws.Range("L10").Value = 3999
$ws.Range("N10").Value = -4339

$ws.Range("H125").Value = 94660.44500000001
$ws.Range("J125").Value = 94660.44500000001
$ws.Range("L125").Value = 94660.44500000001
$ws.Range("N125").Value = -104500.445

$ws.Range("H132").Value = 3234.7222
$ws.Range("I132").Value = 2944.3044
$ws.Range("K132").Value = 8832.913199999999
$ws.Range("M132").Value = -6302.913199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 70448
$ws.Range("J93").Value = 70448
$ws.Range("L93").Value = 70448
$ws.Range("N93").Value = -74192

$ws.Range("H105").Value = 2461.8
$ws.Range("I105").Value = 2072.2222
$ws.Range("K105").Value = 2072.2222
$ws.Range("M105").Value = -325.2222000000002

$ws.Range("H107").Value = 1814.826
$ws.Range("I107").Value = 1732
$ws.Range("J107").Value = 2049.5
$ws.Range("K107").Value = 1732
$ws.Range("L107").Value = 2049.5
$ws.Range("M107").Value = 188
$ws.Range("N107").Value = -5889.5

$ws.Range("H126").Value = 105999
$ws.Range("J126").Value = 105999
$ws.Range("L126").Value = 105999
$ws.Range("N126").Value = -115879

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 59.076923
$ws.Range("I7").Value = 61.333332
$ws.Range("J7").Value = 54
$ws.Range("K7").Value = 61.333332
$ws.Range("L7").Value = 54
$ws.Range("M7").Value = 51.666668
$ws.Range("N7").Value = -280

$ws.Range("H22").Value = 192.5
$ws.Range("I22").Value = 192.5
$ws.Range("K22").Value = 192.5
$ws.Range("M22").Value = 157.5

$ws.Range("H58").Value = 2662.7144
$ws.Range("J58").Value = 2798.261
$ws.Range("L58").Value = 2798.261
$ws.Range("N58").Value = -3204.261

$ws.Range("H97").Value = 105664.25
$ws.Range("J97").Value = 105664.25
$ws.Range("L97").Value = 105664.25
$ws.Range("N97").Value = -107646.25

$ws.Range("H105").Value = 1746.8928
$ws.Range("I105").Value = 1540.1818
$ws.Range("K105").Value = 1540.1818
$ws.Range("M105").Value = 206.8181999999999

$ws.Range("H136").Value = 2662.7144
$ws.Range("J136").Value = 2798.261
$ws.Range("L136").Value = 8394.782999999999
$ws.Range("N136").Value = -13494.783

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 320.0909
$ws.Range("I5").Value = 320.0909
$ws.Range("K5").Value = 960.2727
$ws.Range("M5").Value = -848.2727

$ws.Range("H26").Value = 1054.8334
$ws.Range("J26").Value = 766.3333
$ws.Range("L26").Value = 2298.9999
$ws.Range("N26").Value = -2874.9999

$ws.Range("H29").Value = 3252.875
$ws.Range("J29").Value = 3024.2
$ws.Range("L29").Value = 9072.599999999999
$ws.Range("N29").Value = -9626.599999999999

$ws.Range("H59").Value = 4251.25
$ws.Range("I59").Value = 2005
$ws.Range("J59").Value = 5000
$ws.Range("K59").Value = 6015
$ws.Range("L59").Value = 15000
$ws.Range("M59").Value = -5475
$ws.Range("N59").Value = -16080

$ws.Range("H76").Value = 4013
$ws.Range("I76").Value = 4013
$ws.Range("K76").Value = 12039
$ws.Range("M76").Value = -11656

$ws.Range("H79").Value = 4013
$ws.Range("I79").Value = 4013
$ws.Range("K79").Value = 12039
$ws.Range("M79").Value = -10713

$ws.Range("H129").Value = 780.6
$ws.Range("I129").Value = 456.08334
$ws.Range("K129").Value = 1368.25002
$ws.Range("M129").Value = 3631.74998

$ws.Range("H131").Value = 1189.1936
$ws.Range("I131").Value = 763.9474
$ws.Range("K131").Value = 2291.8422
$ws.Range("M131").Value = 2748.1578

$ws.Range("H135").Value = 320.0909
$ws.Range("I135").Value = 320.0909
$ws.Range("K135").Value = 2880.8181
$ws.Range("M135").Value = -345.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 302
$ws.Range("I2").Value = 65.75
$ws.Range("J2").Value = 774.5
$ws.Range("K2").Value = 65.75
$ws.Range("L2").Value = 774.5
$ws.Range("M2").Value = 47.25
$ws.Range("N2").Value = -1000.5

$ws.Range("H107").Value = 327.34784
$ws.Range("I107").Value = 269.92307
$ws.Range("K107").Value = 269.92307
$ws.Range("M107").Value = 1650.07693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1923.1666
$ws.Range("I16").Value = 1734.3636
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 1734.3636
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -1564.3636
$ws.Range("N16").Value = -4340

$ws.Range("H42").Value = 17527.25
$ws.Range("I42").Value = 10025
$ws.Range("J42").Value = 20028
$ws.Range("K42").Value = 10025
$ws.Range("L42").Value = 20028
$ws.Range("M42").Value = -9462
$ws.Range("N42").Value = -21154

$ws.Range("H46").Value = 3224.3914
$ws.Range("J46").Value = 3748
$ws.Range("L46").Value = 3748
$ws.Range("N46").Value = -4124

$ws.Range("H49").Value = 17527.25
$ws.Range("I49").Value = 10025
$ws.Range("J49").Value = 20028
$ws.Range("K49").Value = 10025
$ws.Range("L49").Value = 20028
$ws.Range("M49").Value = -9878
$ws.Range("N49").Value = -20322

$ws.Range("H61").Value = 1987.8182
$ws.Range("I61").Value = 2606.5
$ws.Range("K61").Value = 2606.5
$ws.Range("M61").Value = -2404.5

$ws.Range("H113").Value = 1987.8182
$ws.Range("I113").Value = 2606.5
$ws.Range("K113").Value = 2606.5
$ws.Range("M113").Value = -436.5

$ws.Range("H122").Value = 14700
$ws.Range("I122").Value = 13707.739
$ws.Range("K122").Value = 41123.217
$ws.Range("M122").Value = -38673.217

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5833.3335
$ws.Range("J81").Value = 10000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -22122

$ws.Range("H84").Value = 5833.3335
$ws.Range("J84").Value = 10000
$ws.Range("L84").Value = 100000
$ws.Range("N84").Value = -110608

$ws.Range("H93").Value = 135000
$ws.Range("J93").Value = 135000
$ws.Range("L93").Value = 135000
$ws.Range("N93").Value = -139992

$ws.Range("H107").Value = 463.64706
$ws.Range("I107").Value = 457.5
$ws.Range("J107").Value = 469.1111
$ws.Range("K107").Value = 1372.5
$ws.Range("L107").Value = 1407.3333
$ws.Range("M107").Value = 547.5
$ws.Range("N107").Value = -5247.3333
